# "Welcome modified for clients"
#
# The original sheet described a "users (psychologist)" / relationship
# table that assumed only one kind of person existed in `users`. The
# edit reshuffles the relationships table (clients/tests/questions/answers),
# renames a couple of field-describing cells for the second ("TABLE
# STRUCTURE") block, and appends three new welcome/briefing questions &
# answers about how the app should work for clients.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Relationships block (rows 4-8) ---------------------------------
# users -> tests -> questions -> answers chain (previously users/clients
# mixed together); rows 7 & 8 no longer carry their own C/D/E relationship
# entries (those rows only keep the right-hand "tables" helper column).
$ws.Range("C4").Value = "users"
$ws.Range("E4").Value = "tests"

$ws.Range("C5").Value = "tests"
$ws.Range("E5").Value = "questions"

$ws.Range("C6").Value = "questions"
$ws.Range("E6").Value = "answers"

$ws.Range("C7:E8").Clear()

# --- TABLE STRUCTURE block ------------------------------------------
# users row: psychologist flag becomes a client flag, plus a new
# psychologist_id column.
$ws.Range("D17").Value = "client(True/False)"
$ws.Range("E17").Value = "psychologist_id"

# answers row: shift the extra columns over (user_id / answer_value).
$ws.Range("F20").Value = "user_id"
$ws.Range("G20").Value = "answer_value"

# --- New welcome / briefing notes for clients (rows 24-26) ----------
$ws.Range("C24").Value = "Is the app generally ment for online use, or it is ment to be used during a live therapy? DURING LIVE THERAPY"
$ws.Range("C25").Value = "Is one app for one psychologist, or one app for many psychologist? ONE APP FOR MANY PSYCHOLOGIST"
$ws.Range("C26").Value = "Should clients log in by themselves, or should they be added by the psycholgist? THEY SHOULD BE ADDED BY THE PSYCHOLOGIST, DURING A LIVE SESSION"

# --- Column E got a touch wider (content is now longer) -------------
$ws.Columns("E").ColumnWidth = 12

# --- Restore view: scrolled down a bit, selection on F22 ------------
$ws.Range("F22").Select()
